# "unit conv and sankey"
# Update the CSET_Set (column G) values on the "TS_Defs Sankey" sheet:
#   - G4: "NRG"  -> "NRG,DEM"
#   - G5: "NRG"  -> formula "=G4" (which evaluates to "NRG,DEM")
# Also move the active selection from C6 to G6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TS_Defs Sankey")

$ws.Range("G4").Value = "NRG,DEM"
$ws.Range("G5").Formula = "=G4"

$ws.Activate()
$ws.Range("G6").Select()
